$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# Update duration for row 22 (C22: 15 -> 20)
$ws.Range("C22").Value = 20

# Row 23: fill in Tache / Duree / Source
$ws.Range("A23").Value = "Explication"
$ws.Range("C23").Value = 30
$ws.Range("E23").Value = "Comment débuter le projet en FORM"

# Row 24: fill in Tache / Duree / Source
$ws.Range("A24").Value = "Code"
$ws.Range("C24").Value = 60
$ws.Range("E24").Value = "Affichage du fond et du joueur"

# Update the selected cell / range in the sheet view
$ws.Range("B25:E25").Select()

$wb.Save()
